# Jogos_da_Semana_FlashScore_2025-03-21.xlsx update
# A new match (Boca Juniors x Real Santander, COLOMBIA - PRIMERA B) is inserted
# as the new row 2, pushing the existing rows down by one. In addition, the
# odds for the "Sportivo Trinidense x Recoleta" match (id IasWL5tH, now row 5)
# were refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (shifts rows 2-5 -> 3-6),
# then strip any formatting the insert may have copied down from the header
# row so the new row matches the plain (unstyled) data rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# New row 2: Boca Juniors vs Real Santander
$ws.Range("A2").Value = "rJpGZjqp"
$ws.Range("B2").Value = "21/03/2025"
$ws.Range("C2").Value = "21:30"
$ws.Range("D2").Value = "COLOMBIA - PRIMERA B"
$ws.Range("E2").Value = "Boca Juniors"
$ws.Range("F2").Value = "Real Santander"
$ws.Range("G2").Value = 2.07
$ws.Range("H2").Value = 2.85
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 2.72
$ws.Range("K2").Value = 1.9
$ws.Range("L2").Value = 4.4
$ws.Range("M2").Value = 1.32
$ws.Range("N2").Value = 2.82
$ws.Range("O2").Value = 1.93
$ws.Range("P2").Value = 1.7
$ws.Range("Q2").Value = 3.1
$ws.Range("R2").Value = 1.27
$ws.Range("S2").Value = 1.47
$ws.Range("T2").Value = 2.32
$ws.Range("U2").Value = 1.7
$ws.Range("V2").Value = 1.93
$ws.Range("W2").Value = 7.1
$ws.Range("X2").Value = 10.25
$ws.Range("Y2").Value = 8.25
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 17
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 8
$ws.Range("AD2").Value = 5.6
$ws.Range("AE2").Value = 12.5
$ws.Range("AF2").Value = 60
$ws.Range("AG2").Value = 450
$ws.Range("AH2").Value = 10.5
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 12.5
$ws.Range("AK2").Value = 70
$ws.Range("AL2").Value = 37
$ws.Range("AM2").Value = 40

# Row 5 is now "IasWL5tH" (Sportivo Trinidense x Recoleta) -- its odds moved
# since the commit and need refreshing.
$ws.Range("G5").Value = 2.25
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 3.25
$ws.Range("J5").Value = 3
$ws.Range("L5").Value = 3.75
$ws.Range("O5").Value = 2.03
$ws.Range("P5").Value = 1.78
$ws.Range("Q5").Value = 3.5
$ws.Range("R5").Value = 1.29
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("W5").Value = 8
$ws.Range("X5").Value = 11
$ws.Range("Y5").Value = 9.5
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 19
$ws.Range("AD5").Value = 6
$ws.Range("AE5").Value = 13
$ws.Range("AH5").Value = 10
$ws.Range("AI5").Value = 15
$ws.Range("AJ5").Value = 12
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 26
$ws.Range("AM5").Value = 34
